# Add the 2024/12/02 data column (CG) to the "合成確率" sheet.
# Mirrors the existing CF column: a date header in row 1 (stored as
# text, like the other date-header cells) plus 52 numeric values in
# rows 2-53, each carrying one of the three existing cell styles
# (s=1 plain Meiryo font, s=2 yellow fill, s=3 light-blue fill).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("合成確率")

# New column CG (#85) needs the same width (12) as the other data columns.
$ws.Columns("CG").ColumnWidth = 11.1666667

# Row 1 header: a literal date string "2024/12/02" (not an Excel date
# serial) -- force text entry via NumberFormat "@" like the neighboring
# header cells, then restyle to match them.
$ws.Range("CG1").NumberFormat = "@"
$ws.Range("CG1").Value = "2024/12/02"
$ws.Range("CF1").Copy()
$ws.Range("CG1").PasteSpecial(-4122)

# Data rows 2-53: value plus the donor cell whose style (font/fill) to
# copy onto it so it reuses the workbook's existing style slots exactly.
$cgData = @(
    @{ Cell = "CG2"; Value = 193.1; Donor = "CF1" },
    @{ Cell = "CG3"; Value = 183.2; Donor = "CF1" },
    @{ Cell = "CG4"; Value = 160.4; Donor = "CF1" },
    @{ Cell = "CG5"; Value = 126.5; Donor = "CF12" },
    @{ Cell = "CG6"; Value = 156.9; Donor = "CF1" },
    @{ Cell = "CG7"; Value = 208.9; Donor = "CF1" },
    @{ Cell = "CG8"; Value = 159.5; Donor = "CF1" },
    @{ Cell = "CG9"; Value = 156.8; Donor = "CF1" },
    @{ Cell = "CG10"; Value = 152.7; Donor = "CF1" },
    @{ Cell = "CG11"; Value = 165.1; Donor = "CF1" },
    @{ Cell = "CG12"; Value = 191.9; Donor = "CF1" },
    @{ Cell = "CG13"; Value = 142.9; Donor = "CF1" },
    @{ Cell = "CG14"; Value = 200.4; Donor = "CF1" },
    @{ Cell = "CG15"; Value = 203.2; Donor = "CF1" },
    @{ Cell = "CG16"; Value = 128.6; Donor = "CF12" },
    @{ Cell = "CG17"; Value = 146.1; Donor = "CF1" },
    @{ Cell = "CG18"; Value = 146.6; Donor = "CF1" },
    @{ Cell = "CG19"; Value = 141.3; Donor = "CF1" },
    @{ Cell = "CG20"; Value = 139.9; Donor = "CF12" },
    @{ Cell = "CG21"; Value = 199.2; Donor = "CF1" },
    @{ Cell = "CG22"; Value = 140.1; Donor = "CF1" },
    @{ Cell = "CG23"; Value = 253.1; Donor = "CF1" },
    @{ Cell = "CG24"; Value = 153.2; Donor = "CF1" },
    @{ Cell = "CG25"; Value = 295; Donor = "CF1" },
    @{ Cell = "CG26"; Value = 149.3; Donor = "CF1" },
    @{ Cell = "CG27"; Value = 150.1; Donor = "CF1" },
    @{ Cell = "CG28"; Value = 126.9; Donor = "CF12" },
    @{ Cell = "CG29"; Value = 126.9; Donor = "CF12" },
    @{ Cell = "CG30"; Value = 160.9; Donor = "CF1" },
    @{ Cell = "CG31"; Value = 265.8; Donor = "CF1" },
    @{ Cell = "CG32"; Value = 189.8; Donor = "CF1" },
    @{ Cell = "CG33"; Value = 173.6; Donor = "CF1" },
    @{ Cell = "CG34"; Value = 139.4; Donor = "CF12" },
    @{ Cell = "CG35"; Value = 173.1; Donor = "CF1" },
    @{ Cell = "CG36"; Value = 156.8; Donor = "CF1" },
    @{ Cell = "CG37"; Value = 150.7; Donor = "CF1" },
    @{ Cell = "CG38"; Value = 129.1; Donor = "CF12" },
    @{ Cell = "CG39"; Value = 111.2; Donor = "CF6" },
    @{ Cell = "CG40"; Value = 170.5; Donor = "CF1" },
    @{ Cell = "CG41"; Value = 170.5; Donor = "CF1" },
    @{ Cell = "CG42"; Value = 243.3; Donor = "CF1" },
    @{ Cell = "CG43"; Value = 135.9; Donor = "CF12" },
    @{ Cell = "CG44"; Value = 217.8; Donor = "CF1" },
    @{ Cell = "CG45"; Value = 127.2; Donor = "CF12" },
    @{ Cell = "CG46"; Value = 238.4; Donor = "CF1" },
    @{ Cell = "CG47"; Value = 163.8; Donor = "CF1" },
    @{ Cell = "CG48"; Value = 219.1; Donor = "CF1" },
    @{ Cell = "CG49"; Value = 125.6; Donor = "CF12" },
    @{ Cell = "CG50"; Value = 440.8; Donor = "CF1" },
    @{ Cell = "CG51"; Value = 160.4; Donor = "CF1" },
    @{ Cell = "CG52"; Value = 137.4; Donor = "CF12" },
    @{ Cell = "CG53"; Value = 144.2; Donor = "CF1" }
)

foreach ($row in $cgData) {
    $ws.Range($row.Cell).Value = $row.Value
    $ws.Range($row.Donor).Copy()
    $ws.Range($row.Cell).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false
